$wb = $excel.ActiveWorkbook

$wsPesos = $wb.Worksheets.Item("Pesos_Locales_Económico")
$wsPesos.Range("B2").Value = 0.06859393436079969
$wsPesos.Range("B3").Value = 0.06859393436079969
$wsPesos.Range("B4").Value = 0.1451001396860344
$wsPesos.Range("B5").Value = 0.06859393436079965
$wsPesos.Range("B6").Value = 0.1451001396860344
$wsPesos.Range("B7").Value = 0.1451001396860344
$wsPesos.Range("B8").Value = 0.01428744057464852
$wsPesos.Range("B9").Value = 0.01428744057464852
$wsPesos.Range("B10").Value = 0.06859393436079965
$wsPesos.Range("B11").Value = 0.0697866565851307
$wsPesos.Range("B12").Value = 0.008867907233478526
$wsPesos.Range("B13").Value = 0.1451001396860344
$wsPesos.Range("B14").Value = 0.009419377695460395
$wsPesos.Range("B15").Value = 0.01428744057464852
$wsPesos.Range("B16").Value = 0.01428744057464852

$wsRanking = $wb.Worksheets.Item("Ranking_Alternativas")
$wsRanking.Range("B2").Value = 0.1127766734448952
$wsRanking.Range("B3").Value = 0.1115091753981816
$wsRanking.Range("B4").Value = 0.1114512740498588
$wsRanking.Range("B5").Value = 0.1101020418150387
$wsRanking.Range("B6").Value = 0.0866053123099499
$wsRanking.Range("B7").Value = 0.07612346363613243
$wsRanking.Range("B8").Value = 0.06619081806160143
$wsRanking.Range("B9").Value = 0.06365593507559968
$wsRanking.Range("B10").Value = 0.06256030094792064
$wsRanking.Range("B11").Value = 0.05013743782618978
$wsRanking.Range("B12").Value = 0.04286718846550173
$wsRanking.Range("B13").Value = 0.03235755559068092
$wsRanking.Range("B14").Value = 0.02914194721354356
$wsRanking.Range("B15").Value = 0.02536492307492816
$wsRanking.Range("B16").Value = 0.01915595308997726

$wsResultados = $wb.Worksheets.Item("Resultados")
$wsResultados.Range("B2").Value = 0.06256030094792064
$wsResultados.Range("B3").Value = 0.06365593507559968
$wsResultados.Range("B4").Value = 0.1101020418150387
$wsResultados.Range("B5").Value = 0.07612346363613243
$wsResultados.Range("B6").Value = 0.1115091753981816
$wsResultados.Range("B7").Value = 0.1127766734448952
$wsResultados.Range("B8").Value = 0.04286718846550173
$wsResultados.Range("B9").Value = 0.01915595308997726
$wsResultados.Range("B10").Value = 0.06619081806160143
$wsResultados.Range("B11").Value = 0.0866053123099499
$wsResultados.Range("B12").Value = 0.05013743782618978
$wsResultados.Range("B13").Value = 0.1114512740498588
$wsResultados.Range("B14").Value = 0.02914194721354356
$wsResultados.Range("B15").Value = 0.03235755559068092
$wsResultados.Range("B16").Value = 0.02536492307492816

$wsMatriz = $wb.Worksheets.Item("Matriz_Económico")
$wsMatriz.Range("D2").Value = 0.3333333333333333
$wsMatriz.Range("F2").Value = 0.3333333333333333
$wsMatriz.Range("G2").Value = 0.3333333333333333
$wsMatriz.Range("M2").Value = 0.3333333333333333
$wsMatriz.Range("N2").Value = 7
$wsMatriz.Range("D3").Value = 0.3333333333333333
$wsMatriz.Range("F3").Value = 0.3333333333333333
$wsMatriz.Range("G3").Value = 0.3333333333333333
$wsMatriz.Range("M3").Value = 0.3333333333333333
$wsMatriz.Range("N3").Value = 7
$wsMatriz.Range("B4").Value = 3
$wsMatriz.Range("C4").Value = 3
$wsMatriz.Range("E4").Value = 3
$wsMatriz.Range("J4").Value = 3
$wsMatriz.Range("K4").Value = 3
$wsMatriz.Range("D5").Value = 0.3333333333333333
$wsMatriz.Range("F5").Value = 0.3333333333333333
$wsMatriz.Range("G5").Value = 0.3333333333333333
$wsMatriz.Range("M5").Value = 0.3333333333333333
$wsMatriz.Range("N5").Value = 7
$wsMatriz.Range("B6").Value = 3
$wsMatriz.Range("C6").Value = 3
$wsMatriz.Range("E6").Value = 3
$wsMatriz.Range("J6").Value = 3
$wsMatriz.Range("K6").Value = 3
$wsMatriz.Range("B7").Value = 3
$wsMatriz.Range("C7").Value = 3
$wsMatriz.Range("E7").Value = 3
$wsMatriz.Range("J7").Value = 3
$wsMatriz.Range("K7").Value = 3
$wsMatriz.Range("D10").Value = 0.3333333333333333
$wsMatriz.Range("F10").Value = 0.3333333333333333
$wsMatriz.Range("G10").Value = 0.3333333333333333
$wsMatriz.Range("M10").Value = 0.3333333333333333
$wsMatriz.Range("N10").Value = 7
$wsMatriz.Range("D11").Value = 0.3333333333333333
$wsMatriz.Range("F11").Value = 0.3333333333333333
$wsMatriz.Range("G11").Value = 0.3333333333333333
$wsMatriz.Range("M11").Value = 0.3333333333333333
$wsMatriz.Range("B13").Value = 3
$wsMatriz.Range("C13").Value = 3
$wsMatriz.Range("E13").Value = 3
$wsMatriz.Range("J13").Value = 3
$wsMatriz.Range("K13").Value = 3
$wsMatriz.Range("B14").Value = 0.1428571428571428
$wsMatriz.Range("C14").Value = 0.1428571428571428
$wsMatriz.Range("E14").Value = 0.1428571428571428
$wsMatriz.Range("J14").Value = 0.1428571428571428

Write-Output "OK"
